$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

$ws.Cells.Item($row, 1).Value = 33
$ws.Cells.Item($row, 2).Value = "india"
$ws.Cells.Item($row, 3).Value = "isl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45235.64583333334
$ws.Cells.Item($row, 6).Value = "Chennaiyin"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Goa"
$ws.Cells.Item($row, 9).Value = 3
$ws.Cells.Item($row, 10).Value = 2.41
$ws.Cells.Item($row, 11).Value = "29/10/2023 15:42"
$ws.Cells.Item($row, 12).Value = 3.05
$ws.Cells.Item($row, 13).Value = "05/11/2023 15:29"
$ws.Cells.Item($row, 14).Value = 3.52
$ws.Cells.Item($row, 15).Value = "29/10/2023 15:42"
$ws.Cells.Item($row, 16).Value = 3.97
$ws.Cells.Item($row, 17).Value = "05/11/2023 15:29"
$ws.Cells.Item($row, 18).Value = 2.82
$ws.Cells.Item($row, 19).Value = "29/10/2023 15:42"
$ws.Cells.Item($row, 20).Value = 2.15
$ws.Cells.Item($row, 21).Value = "05/11/2023 15:29"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/india/isl/chennaiyin-fc-fc-goa/KUlAzUpM/"

# Copy the number/cell formats from the row above so the new row matches
# the existing style table (bordered/bold index col, datetime serial col)
# instead of Excel synthesizing brand-new style entries.
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null

$ws.Range("E33").Copy() | Out-Null
$ws.Range("E34").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
